$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$new1 = $wb.Worksheets.Add($null, $ws)
$ws.Range("A1:E6").Copy($new1.Range("A1"))
$new1.Hyperlinks.Add($new1.Range("C2"), "mailto:abc@gamail.com")
